$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("interactions")
$ws.Range("A1").Value = "TEST"
